$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp string in cell A1
$ws.Range("A1").Value = "Datos actualizados a 15 de Septiembre de 2020 a las 15:29"

# Update country statistics cells per the diff
$ws.Range("B4").Value = 6751048
$ws.Range("C4").Value = 1759
$ws.Range("D4").Value = 4028756
$ws.Range("E4").Value = 2523144
$ws.Range("G4").Value = 148
$ws.Range("H4").Value = 199148
$ws.Range("B5").Value = 4933188
$ws.Range("C5").Value = 6274
$ws.Range("E5").Value = 992956
$ws.Range("G5").Value = 25
$ws.Range("H5").Value = 80833
$ws.Range("D13").Value = 438883
$ws.Range("E13").Value = 114853
$ws.Range("G13").Value = 43
$ws.Range("H13").Value = 11710
$ws.Range("B19").Value = 326930
$ws.Range("C19").Value = 672
$ws.Range("D19").Value = 305022
$ws.Range("E19").Value = 17570
$ws.Range("G19").Value = 33
$ws.Range("H19").Value = 4338
$ws.Range("B31").Value = 122214
$ws.Range("C31").Value = 239
$ws.Range("D31").Value = 119144
$ws.Range("E31").Value = 2862
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 208
$ws.Range("B38").Value = 96301
$ws.Range("C38").Value = 829
$ws.Range("D38").Value = 86219
$ws.Range("E38").Value = 9514
$ws.Range("G38").Value = 5
$ws.Range("H38").Value = 568
$ws.Range("B42").Value = 87345
$ws.Range("G42").Value = 4
$ws.Range("H42").Value = 5851
$ws.Range("B44").Value = 84778
$ws.Range("C44").Value = 1379
$ws.Range("G44").Value = 2
$ws.Range("H44").Value = 6258
$ws.Range("B46").Value = 80940
$ws.Range("C46").Value = 674
$ws.Range("D46").Value = 70635
$ws.Range("E46").Value = 9904
$ws.Range("G46").Value = 2
$ws.Range("H46").Value = 401
$ws.Range("B49").Value = 74552
$ws.Range("C49").Value = 192
$ws.Range("D49").Value = 72661
$ws.Range("E49").Value = 1130
$ws.Range("G49").Value = 5
$ws.Range("H49").Value = 761
$ws.Range("B60").Value = 48429
$ws.Range("C60").Value = 593
$ws.Range("D60").Value = 44942
$ws.Range("E60").Value = 3085
$ws.Range("G60").Value = 6
$ws.Range("H60").Value = 402
$ws.Range("B67").Value = 38517
$ws.Range("C67").Value = 114
$ws.Range("D67").Value = 35998
$ws.Range("E67").Value = 1953
$ws.Range("G67").Value = 2
$ws.Range("H67").Value = 566
$ws.Range("B71").Value = 32511
$ws.Range("C71").Value = 74
$ws.Range("D71").Value = 31313
$ws.Range("E71").Value = 463
$ws.Range("G71").Value = 2
$ws.Range("H71").Value = 735
$ws.Range("B72").Value = 32250
$ws.Range("C72").Value = 888
$ws.Range("D72").Value = 21804
$ws.Range("E72").Value = 10217
$ws.Range("G72").Value = 3
$ws.Range("H72").Value = 229
$ws.Range("B85").Value = 15925
$ws.Range("C85").Value = 98
$ws.Range("D85").Value = 13418
$ws.Range("E85").Value = 1846
$ws.Range("G85").Value = 9
$ws.Range("H85").Value = 661
$ws.Range("E102").Value = 886
$ws.Range("G102").Value = 2
$ws.Range("H102").Value = 339
$ws.Range("B136").Value = 3266
$ws.Range("C136").Value = 4
$ws.Range("D136").Value = 3016
$ws.Range("E136").Value = 237
$ws.Range("B150").Value = 2174
$ws.Range("C150").Value = 6
$ws.Range("D150").Value = 2102
$ws.Range("E150").Value = 62
$ws.Range("E203").Value = 10
$ws.Range("H203").Value = 1
